$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the values in row 4 (columns E through K)
$ws.Range("E4").Value = 2703
$ws.Range("F4").Value = 2722
$ws.Range("G4").Value = 2702
$ws.Range("H4").Value = 2697
$ws.Range("I4").Value = 2830
$ws.Range("J4").Value = 2854
$ws.Range("K4").Value = 2879

# Update the active selection to A3
$ws.Range("A3").Select()
